# The commit swaps the contents of ppt/theme/theme1.xml (the default
# "Office Theme" colour scheme) and ppt/theme/theme2.xml (the "Integral" /
# "Red Violet" colour scheme) while leaving the font scheme and format
# scheme (identical between the two themes) untouched.
#
# theme2.xml is the theme actually bound to the slide master / presentation
# (and so is the one reachable through the PowerPoint object model's
# ThemeColorScheme). We drive its 12 colour scheme slots to the values
# that used to live in theme1.xml ("Office Theme" colours), matching the
# post-edit content of theme2.xml in the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Order of ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# Values below are the "Office" colour scheme (hex -> COM RGB long =
# R + G*256 + B*65536) that was previously stored in theme1.xml.

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
